$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5569
$ws.Range("E2").Value = 119
$ws.Range("F2").Value = 119
$ws.Range("G2").Value = 726
$ws.Range("H2").Value = 503
$ws.Range("I2").Value = 508
$ws.Range("J2").Value = -4
$ws.Range("K2").Value = 12729
$ws.Range("L2").Value = 2300
$ws.Range("M2").Value = 10429
$ws.Range("N2").Value = 10270
$ws.Range("O2").Value = 158
$ws.Range("P2").Value = 76
$ws.Range("Q2").Value = 227
$ws.Range("R2").Value = -901
$ws.Range("S2").Value = -12
$ws.Range("T2").Value = 286
$ws.Range("U2").Value = -59
$ws.Range("V2").Value = 566
$ws.Range("W2").Value = 2.13
$ws.Range("X2").Value = 9.039999999999999
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 4.11
$ws.Range("AA2").Value = 22.05
$ws.Range("AB2").Value = 13774.5
$ws.Range("AC2").Value = 3236
$ws.Range("AD2").Value = 26.15
$ws.Range("AE2").Value = 66969
$ws.Range("AF2").Value = 1.26
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 15691955

# Row 3
$ws.Range("D3").Value = 6446
$ws.Range("E3").Value = -543
$ws.Range("F3").Value = -543
$ws.Range("G3").Value = 2290
$ws.Range("H3").Value = 1652
$ws.Range("I3").Value = 1640
$ws.Range("J3").Value = 11
$ws.Range("K3").Value = 18333
$ws.Range("L3").Value = 3308
$ws.Range("M3").Value = 15025
$ws.Range("N3").Value = 14270
$ws.Range("O3").Value = 755
$ws.Range("P3").Value = 98
$ws.Range("Q3").Value = -1075
$ws.Range("R3").Value = -1175
$ws.Range("S3").Value = 2771
$ws.Range("T3").Value = 698
$ws.Range("U3").Value = -1772
$ws.Range("V3").Value = 625
$ws.Range("W3").Value = -8.43
$ws.Range("X3").Value = 25.62
$ws.Range("Y3").Value = 13.37
$ws.Range("Z3").Value = 10.63
$ws.Range("AA3").Value = 22.02
$ws.Range("AB3").Value = 15216.36
$ws.Range("AC3").Value = 8836
$ws.Range("AD3").Value = 6.44
$ws.Range("AE3").Value = 74287
$ws.Range("AF3").Value = 0.77
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 19565025

# Row 4
$ws.Range("D4").Value = 8564
$ws.Range("E4").Value = 264
$ws.Range("F4").Value = 264
$ws.Range("G4").Value = 260
$ws.Range("H4").Value = 71
$ws.Range("I4").Value = 125
$ws.Range("J4").Value = -54
$ws.Range("K4").Value = 18382
$ws.Range("L4").Value = 3277
$ws.Range("M4").Value = 15104
$ws.Range("N4").Value = 14628
$ws.Range("O4").Value = 476
$ws.Range("P4").Value = 98
$ws.Range("Q4").Value = 834
$ws.Range("R4").Value = 367
$ws.Range("S4").Value = -792
$ws.Range("T4").Value = 356
$ws.Range("U4").Value = 478
$ws.Range("V4").Value = 231
$ws.Range("W4").Value = 3.08
$ws.Range("X4").Value = 0.83
$ws.Range("Y4").Value = 0.87
$ws.Range("Z4").Value = 0.39
$ws.Range("AA4").Value = 21.7
$ws.Range("AB4").Value = 15336.17
$ws.Range("AC4").Value = 639
$ws.Range("AD4").Value = 78.66
$ws.Range("AE4").Value = 76153
$ws.Range("AF4").Value = 0.66
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 19565025

# Row 5
$ws.Range("D5").Value = 9091
$ws.Range("E5").Value = 347
$ws.Range("F5").Value = 347
$ws.Range("G5").Value = 544
$ws.Range("H5").Value = 87
$ws.Range("I5").Value = 166
$ws.Range("J5").Value = -79
$ws.Range("K5").Value = 21176
$ws.Range("L5").Value = 4597
$ws.Range("M5").Value = 16579
$ws.Range("N5").Value = 15066
$ws.Range("O5").Value = 1513
$ws.Range("P5").Value = 98
$ws.Range("Q5").Value = -435
$ws.Range("R5").Value = -857
$ws.Range("S5").Value = 1962
$ws.Range("T5").Value = 840
$ws.Range("U5").Value = -1275
$ws.Range("V5").Value = 978
$ws.Range("W5").Value = 3.82
$ws.Range("X5").Value = 0.96
$ws.Range("Y5").Value = 1.12
$ws.Range("Z5").Value = 0.44
$ws.Range("AA5").Value = 27.73
$ws.Range("AB5").Value = 15497.19
$ws.Range("AC5").Value = 849
$ws.Range("AD5").Value = 80.61
$ws.Range("AE5").Value = 78434
$ws.Range("AF5").Value = 0.87
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 19565025

# Row 6
$ws.Range("D6").Value = 12646
$ws.Range("E6").Value = 686
$ws.Range("F6").Value = 686
$ws.Range("G6").Value = 1538
$ws.Range("H6").Value = 1047
$ws.Range("I6").Value = 894
$ws.Range("K6").Value = 24947
$ws.Range("L6").Value = 6746
$ws.Range("M6").Value = 18200
$ws.Range("N6").Value = 15823
$ws.Range("P6").Value = 98
$ws.Range("Q6").Value = 381
$ws.Range("R6").Value = 234
$ws.Range("S6").Value = -443
$ws.Range("T6").Value = 740
$ws.Range("U6").Value = -359
$ws.Range("V6").Value = 1409
$ws.Range("W6").Value = 5.43
$ws.Range("X6").Value = 8.279999999999999
$ws.Range("Y6").Value = 5.79
$ws.Range("Z6").Value = 4.54
$ws.Range("AA6").Value = 37.07
$ws.Range("AB6").Value = 16599.21
$ws.Range("AC6").Value = 4567
$ws.Range("AD6").Value = 12.55
$ws.Range("AE6").Value = 83901
$ws.Range("AF6").Value = 0.68
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 19565025

# Row 7
$ws.Range("D7").Value = 14745
$ws.Range("E7").Value = 942
$ws.Range("G7").Value = 1334
$ws.Range("H7").Value = 840
$ws.Range("I7").Value = 695
$ws.Range("K7").Value = 26130
$ws.Range("L7").Value = 7122
$ws.Range("M7").Value = 19007
$ws.Range("N7").Value = 16550
$ws.Range("P7").Value = 99
$ws.Range("Q7").Value = 1899
$ws.Range("R7").Value = -1007
$ws.Range("S7").Value = -324
$ws.Range("T7").Value = 857
$ws.Range("U7").Value = 1166
$ws.Range("W7").Value = 6.39
$ws.Range("X7").Value = 5.7
$ws.Range("Y7").Value = 4.29
$ws.Range("Z7").Value = 3.29
$ws.Range("AA7").Value = 37.47
$ws.Range("AC7").Value = 3554
$ws.Range("AD7").Value = 21.25
$ws.Range("AE7").Value = 88558
$ws.Range("AF7").Value = 0.85
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0

# Row 8
$ws.Range("D8").Value = 16100
$ws.Range("E8").Value = 1208
$ws.Range("G8").Value = 1615
$ws.Range("H8").Value = 1121
$ws.Range("I8").Value = 957
$ws.Range("K8").Value = 27932
$ws.Range("L8").Value = 7824
$ws.Range("M8").Value = 20109
$ws.Range("N8").Value = 17429
$ws.Range("P8").Value = 99
$ws.Range("Q8").Value = 1942
$ws.Range("R8").Value = -1136
$ws.Range("S8").Value = 31
$ws.Range("T8").Value = 891
$ws.Range("U8").Value = 1136
$ws.Range("W8").Value = 7.5
$ws.Range("X8").Value = 6.96
$ws.Range("Y8").Value = 5.63
$ws.Range("Z8").Value = 4.15
$ws.Range("AA8").Value = 38.91
$ws.Range("AC8").Value = 4891
$ws.Range("AD8").Value = 15.44
$ws.Range("AE8").Value = 93257
$ws.Range("AF8").Value = 0.8100000000000001
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0

# Row 9
$ws.Range("D9").Value = 17569
$ws.Range("E9").Value = 1360
$ws.Range("G9").Value = 1797
$ws.Range("H9").Value = 1258
$ws.Range("I9").Value = 1082
$ws.Range("K9").Value = 30106
$ws.Range("L9").Value = 8746
$ws.Range("M9").Value = 21358
$ws.Range("N9").Value = 18361
$ws.Range("P9").Value = 99
$ws.Range("Q9").Value = 2074
$ws.Range("R9").Value = -1039
$ws.Range("S9").Value = -26
$ws.Range("T9").Value = 897
$ws.Range("U9").Value = 1446
$ws.Range("W9").Value = 7.74
$ws.Range("X9").Value = 7.16
$ws.Range("Y9").Value = 6.05
$ws.Range("Z9").Value = 4.33
$ws.Range("AA9").Value = 40.95
$ws.Range("AC9").Value = 5532
$ws.Range("AD9").Value = 13.65
$ws.Range("AE9").Value = 98244
$ws.Range("AF9").Value = 0.77
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0

# Remove cells that no longer exist in the target (ClearContents removes the <c> entirely)
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("AI9").ClearContents()
